$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.49"
$ws.Range("E2").Value = "'-0.53%"
$ws.Range("D3").Value = "'27.11"
$ws.Range("E3").Value = "'3.83%"
$ws.Range("D4").Value = "'5.162"
$ws.Range("E4").Value = "'1.58%"
$ws.Range("D5").Value = "'0.05651"
$ws.Range("E5").Value = "'0.83%"
$ws.Range("D6").Value = "'6.471"
$ws.Range("E6").Value = "'-0.13%"
$ws.Range("D7").Value = "'0.8198"
$ws.Range("E7").Value = "'0.79%"
$ws.Range("D8").Value = "'0.8458"
$ws.Range("E8").Value = "'0.22%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1329"
$ws.Range("E9").Value = "'-1.24%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.06937"
$ws.Range("E10").Value = "'-1.02%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.02930"
$ws.Range("E11").Value = "'4.80%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09404"
$ws.Range("E12").Value = "'0.14%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001519"
$ws.Range("E13").Value = "'-0.42%"
$ws.Range("B14").Value = "CoinExToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D14").Value = "'0.04307"
$ws.Range("E14").Value = "'-8.04%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005978"
$ws.Range("E15").Value = "'-0.87%"
$ws.Range("D16").Value = "'0.006147"
$ws.Range("E16").Value = "'0.14%"
$ws.Range("D17").Value = "'3.510"
$ws.Range("E17").Value = "'-1.31%"
$ws.Range("D18").Value = "'3.003"
$ws.Range("E18").Value = "'-1.04%"
$ws.Range("E19").Value = "'9.12%"
$ws.Range("D20").Value = "'0.3114"
$ws.Range("E20").Value = "'-2.33%"
$ws.Range("D21").Value = "'0.03120"
$ws.Range("E21").Value = "'-2.98%"
$ws.Range("D22").Value = "'0.1258"
$ws.Range("E22").Value = "'-4.66%"
$ws.Range("D23").Value = "'2.127"
$ws.Range("E23").Value = "'-43.38%"
$ws.Range("D24").Value = "'0.1373"
$ws.Range("E24").Value = "'-0.12%"
$ws.Range("E25").Value = "'-1.77%"
$ws.Range("D26").Value = "'0.004473"
$ws.Range("E26").Value = "'-2.81%"
$ws.Range("D27").Value = "'0.00009797"
$ws.Range("D28").Value = "'0.00007255"
$ws.Range("E28").Value = "'-47.80%"
$ws.Range("E40").Value = "'-0.28%"
$ws.Range("D41").Value = "'0.006088"
$ws.Range("E41").Value = "'-1.03%"
$ws.Range("D42").Value = "'0.1053"
$ws.Range("E42").Value = "'-0.12%"
$ws.Range("D43").Value = "'0.002299"
$ws.Range("E43").Value = "'-8.04%"
$ws.Range("E44").Value = "'-6.45%"
$ws.Range("D45").Value = "'0.00005362"
$ws.Range("E45").Value = "'1.32%"
$ws.Range("E46").Value = "'-0.04%"
$ws.Range("E47").Value = "'-23.98%"
$ws.Range("D48").Value = "'0.002650"
$ws.Range("E48").Value = "'28.76%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.04%"
